$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '27.441.41'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -2.57%  '

# Row 3
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.740.94'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -3.41%  '

# Row 4
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.003'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.15%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '323.36'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -4.50%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.9988'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.44%  '

# Row 7
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4259'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -9.35%  '

# Row 8
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3615'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -5.16%  '

# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '45.02'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -1.12%  '

# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.131'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -1.11%  '

# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07420'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -3.47%  '

# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.9989'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -0.50%  '

# Row 13
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '21.80'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -2.78%  '

# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.064'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -5.28%  '

# Row 15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.179'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -2.58%  '

# Row 16
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '1.735.29'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -3.60%  '

# Row 17
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.00001063'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -3.21%  '

# Row 18
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '87.80'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +6.98%  '

# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.05968'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -11.70%  '

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.9997'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.10%  '

# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '16.91'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -3.30%  '

# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.074'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -5.82%  '

# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.5218'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -5.45%  '

# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '27.486.51'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -2.35%  '

# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '11.30'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -5.66%  '

# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.402'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.35%  '

# Row 27
$ws.Range('B27').Value = 'EthereumClassic'
$ws.Range('C27').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '20.13'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -3.67%  '

# Row 28
$ws.Range('B28').Value = 'LidoDAOToken'
$ws.Range('C28').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.371'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -1.53%  '

# Row 29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '149.53'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -1.32%  '

# Row 30
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.934.70'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -3.66%  '

# Row 31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.271'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.05%  '

# Row 32
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '126.27'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -5.74%  '

# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.701'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -8.99%  '

# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.605'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -5.75%  '

# Row 35
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.09009'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -7.80%  '

# Row 36
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '12.50'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +2.51%  '

# Row 37
$ws.Range('E37').Value = '  -3.56%  '

# Row 38
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.02281'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -4.41%  '

# Row 39
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.06143'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -3.62%  '

# Row 40
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.6444'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -3.81%  '

# Row 41
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '5.022'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -4.22%  '

# Row 42
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.183'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -3.42%  '

# Row 43
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.424'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -4.64%  '

# Row 44
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.9999'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.14%  '

# Row 45
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '7.820'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -3.37%  '

# Row 46
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '13.47'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -4.51%  '

# Row 47
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.743'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -3.35%  '

# Row 48
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.5882'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -4.59%  '

# Row 49
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '125.32'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -3.12%  '

# Row 50
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.943'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -5.55%  '

# Row 51
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.06829'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -4.15%  '
